$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44181
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 12000
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 3000
$ws.Range("O2").Value = "Provincia de Chacabuco"
$ws.Range("P2").Value = 30
$ws.Range("D3").Value = 44245
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 9000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 30
$ws.Range("D4").Value = 44245
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2500
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 25
$ws.Range("D5").Value = 44232
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 16000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("O5").Value = "Provincia de Chacabuco"
$ws.Range("P5").Value = 30
$ws.Range("D6").Value = 44159
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 7000
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 3000
$ws.Range("O6").Value = "Provincia de Chacabuco"
$ws.Range("P6").Value = 30
$ws.Range("D7").Value = 44188
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 12000
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("O7").Value = "Provincia de Chacabuco"
$ws.Range("P7").Value = 30
$ws.Range("D8").Value = 44189
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 16000
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("O8").Value = "Provincia de Chacabuco"
$ws.Range("P8").Value = 30
$ws.Range("D9").Value = 44166
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 7000
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 3000
$ws.Range("O9").Value = "Provincia de Chacabuco"
$ws.Range("P9").Value = 30
$ws.Range("D10").Value = 44231
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 12000
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("O10").Value = "Provincia de Chacabuco"
$ws.Range("P10").Value = 30
$ws.Range("D11").Value = 44187
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 12000
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 3000
$ws.Range("O11").Value = "Provincia de Chacabuco"
$ws.Range("P11").Value = 30
$ws.Range("D12").Value = 44214
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 7000
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("O12").Value = "Provincia de Chacabuco"
$ws.Range("P12").Value = 30
$ws.Range("D13").Value = 44215
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 16000
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 3000
$ws.Range("O13").Value = "Provincia de Chacabuco"
$ws.Range("P13").Value = 30
$ws.Range("D14").Value = 44204
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 7000
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = 3000
$ws.Range("O14").Value = "Provincia de Chacabuco"
$ws.Range("P14").Value = 30
$ws.Range("D15").Value = 44168
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 7000
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 3000
$ws.Range("O15").Value = "Provincia de Chacabuco"
$ws.Range("P15").Value = 30
$ws.Range("D16").Value = 44161
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("O16").Value = "Provincia de Chacabuco"
$ws.Range("P16").Value = 30
$ws.Range("D17").Value = 44160
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 7000
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 3000
$ws.Range("O17").Value = "Provincia de Chacabuco"
$ws.Range("P17").Value = 30
$ws.Range("D18").Value = 44210
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 8800
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = 2750
$ws.Range("O18").Value = "Provincia de Chacabuco"
$ws.Range("P18").Value = 28
$ws.Range("D19").Value = 44230
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 16000
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("O19").Value = "Provincia de Chacabuco"
$ws.Range("P19").Value = 30
$ws.Range("D20").Value = 44186
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = 3000
$ws.Range("O20").Value = "Provincia de Chacabuco"
$ws.Range("P20").Value = 30
$ws.Range("D21").Value = 44167
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 7000
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = 3000
$ws.Range("O21").Value = "Provincia de Chacabuco"
$ws.Range("P21").Value = 30
$ws.Range("D22").Value = 44209
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 7000
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 2750
$ws.Range("O22").Value = "Provincia de Chacabuco"
$ws.Range("P22").Value = 28
$ws.Range("D23").Value = 44229
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 16000
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = 3000
$ws.Range("O23").Value = "Provincia de Chacabuco"
$ws.Range("P23").Value = 30
$ws.Range("D24").Value = 44162
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 7000
$ws.Range("K24").Value = 3000
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = 3000
$ws.Range("O24").Value = "Provincia de Chacabuco"
$ws.Range("P24").Value = 30
